# Apply the "StoreIt" pricing-page rewrite described by the target diff.
#
# Strategy:
#  - Word's Find/Replace (the Replacement.Text path) runs the replacement
#    text through AutoFormat/AutoCorrect (e.g. straight quotes become
#    curly/smart quotes), which the target XML does NOT want. So
#    Find.Execute() is used only to *locate* text; mutation is always a
#    direct Range.Text assignment (plain replace, no autoformat) or
#    Range.InsertXML (when a run needs to carry a <w:proofErr/> split or
#    other run-level markup - InsertXML is the documented way to inject
#    that exact markup. It REPLACES the contents of the range it is
#    called on, so it is always called on a whole-paragraph range to
#    swap a paragraph's runs wholesale).
#  - After Find.Execute(), the returned Range does not automatically span
#    the whole paragraph (and Paragraphs.First/.Next on it are not
#    reliable ranges in this host), so Range.Expand(4) (wdParagraph) is
#    used to grow the hit to the enclosing paragraph, and the paragraph's
#    1-based index is recovered by counting paragraphs from the start of
#    the document up to the range - giving a plain $d.Paragraphs(i) we
#    can keep re-using (including i+1, i+2, ... for freshly inserted
#    paragraphs).
#  - Edits are applied from the BOTTOM of the document upward so that
#    inserting/expanding paragraphs never invalidates paragraph indices
#    still needed later in the script.

$d = $word.ActiveDocument
$wdParagraph = 4

function Get-ParaIndex($rng) {
    return $d.Range(0, $rng.End).Paragraphs.Count
}

function New-WordPkgXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Replace the whole paragraph located by a unique text anchor with the
# given inner <w:p>...</w:p>+ markup (one or more paragraphs).
function Replace-ParagraphByAnchor([string]$anchor, [string]$newParasXml) {
    $r = $d.Content
    $null = $r.Find.Execute($anchor)
    $r.Expand($wdParagraph)
    $r.InsertXML((New-WordPkgXml($newParasXml)))
}

# ---------------------------------------------------------------------
# 1. Footer: "&copy; Your Company Name 2024" -> "<p>&copy; StoreIt 2024</p>"
# ---------------------------------------------------------------------
$xmlFooter = '<w:p><w:r><w:t xml:space="preserve">    &lt;p&gt;&amp;copy; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>StoreIt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 2024&lt;/p&gt;</w:t></w:r></w:p>'
Replace-ParagraphByAnchor "&copy; Your Company Name 2024" $xmlFooter

# ---------------------------------------------------------------------
# 2. Insert the "Basic" (popular) and "Pro" pricing-plan <div> blocks right
#    after the Free plan's closing </div>, before the pricing-container's
#    closing </div>.
# ---------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("Sign Up</a>")
$r.Expand($wdParagraph)
$signUpIdx = Get-ParaIndex($r)
$freePlanDivEnd = $d.Paragraphs($signUpIdx + 1)
$freePlanDivEnd.Range.InsertParagraphAfter()
$newSpot = $d.Paragraphs($signUpIdx + 2)
$xmlPlans = '<w:p><w:r><w:t xml:space="preserve">      &lt;div class="pricing-plan popular"&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        &lt;h3&gt;Basic&lt;/h3&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        &lt;p&gt;25GB Storage&lt;/p&gt;</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">        &lt;p&gt;File Versioning &amp; Syncing&lt;/p&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        &lt;p&gt;Mobile &amp; Desktop Apps&lt;/p&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>="#"&gt;Sign Up&lt;/a&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        &lt;span class="popular-label"&gt;Popular&lt;/span&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      &lt;/div&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      &lt;div class="pricing-plan"&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        &lt;h3&gt;Pro&lt;/h3&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        &lt;p&gt;100GB Storage&lt;/p&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        &lt;p&gt;Advanced Security Features&lt;/p&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        &lt;p&gt;24/7 Customer Support&lt;/p&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>="#"&gt;Sign Up&lt;/a&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      &lt;/div&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    &lt;/div&gt;</w:t></w:r></w:p>'
$newSpot.Range.InsertXML((New-WordPkgXml($xmlPlans)))

# ---------------------------------------------------------------------
# 3. Free plan: "<p>Basic features...</p>" -> "<p>5GB Storage</p>" plus a
#    new "<p>Basic File Access</p>" paragraph.
# ---------------------------------------------------------------------
$xmlFree = '<w:p><w:r><w:t xml:space="preserve">        &lt;p&gt;5GB Storage&lt;/p&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        &lt;p&gt;Basic File Access&lt;/p&gt;</w:t></w:r></w:p>'
Replace-ParagraphByAnchor "Basic features..." $xmlFree

# ---------------------------------------------------------------------
# 4. "<h2>Pricing Plans</h2>" -> "<h2>Choose the Perfect Plan</h2>"
# ---------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("Pricing Plans")
$r.Text = "Choose the Perfect Plan"

# ---------------------------------------------------------------------
# 5. Header: "<h1>Your Company Name</h1>" -> "<h1>StoreIt - Your Cloud
#    Storage Solution</h1>", and "<nav>...</nav>" expands into an open
#    <nav> tag, three links, and a closing </nav> tag (each its own
#    paragraph).
# ---------------------------------------------------------------------
$xmlH1 = '<w:p><w:r><w:t xml:space="preserve">    &lt;h1&gt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>StoreIt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - Your Cloud Storage Solution&lt;/h1&gt;</w:t></w:r></w:p>'
Replace-ParagraphByAnchor "Your Company Name</h1>" $xmlH1

$xmlNav = '<w:p><w:r><w:t xml:space="preserve">    &lt;nav&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>="#"&gt;Features&lt;/a&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>="#"&gt;Security&lt;/a&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">      &lt;a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>href</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>="#"&gt;Contact&lt;/a&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    &lt;/nav&gt;</w:t></w:r></w:p>'
Replace-ParagraphByAnchor "<nav>...</nav>" $xmlNav

# ---------------------------------------------------------------------
# 6. "<title>Your Pricing | Product Name</title>" -> "<title>Cloud Storage
#    Plans | StoreIt</title>"
# ---------------------------------------------------------------------
$xmlTitle = '<w:p><w:r><w:t xml:space="preserve">  &lt;title&gt;Cloud Storage Plans | </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>StoreIt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&lt;/title&gt;</w:t></w:r></w:p>'
Replace-ParagraphByAnchor "Your Pricing | Product Name" $xmlTitle

# ---------------------------------------------------------------------
# 7. '<html lang="' -> '<lang="' (only the first run's text changes; the
#    "en" proofErr-wrapped run and the closing '">' run are untouched).
# ---------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute('<html lang="')
$r.Text = '<lang="'

Write-Host "Edit complete."
